$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("plateK10rep1_20200429_122048_065_WellA06_PointA06_0007_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0052.h5", "second_pass"),
    @("plateK10rep1_20200429_122048_065_WellG06_PointG06_0005_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0662.h5", "second_pass"),
    @("plateK10rep1_20200429_122048_065_WellH02_PointH02_0003_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0813.h5", "ok"),
    @("plateK11rep1_20200429_140316_208_WellD05_PointD05_0007_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0394.h5", "second_pass"),
    @("plateK11rep1_20200429_140316_208_WellE05_PointE05_0004_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0454.h5", "second_pass"),
    @("plateK11rep1_20200429_140316_208_WellH06_PointH06_0006_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0780.h5", "second_pass"),
    @("plateK12rep1_20200430_155932_313_WellC08_PointC08_0000_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0279.h5", "second_pass"),
    @("plateK12rep1_20200430_155932_313_WellD09_PointD09_0003_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0354.h5", "second_pass"),
    @("plateK12rep1_20200430_155932_313_WellF04_PointF04_0002_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0614.h5", "second_pass"),
    @("plateK13rep1_20200430_175056_461_WellA10_PointA10_0007_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0088.h5", "second_pass"),
    @("plateK13rep1_20200430_175056_461_WellE04_PointE04_0000_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0108.h5", "ok"),
    @("plateK13rep1_20200430_175056_461_WellF02_PointF02_0008_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0161.h5", "second_pass"),
    @("plateK13rep1_20200430_175056_461_WellF08_PointF08_0008_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0215.h5", "second_pass"),
    @("plateK14rep1_20200430_194338_941_WellE09_PointE09_0002_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0506.h5", "second_pass"),
    @("plateK14rep1_20200430_194338_941_WellF02_PointF02_0003_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0633.h5", "second_pass"),
    @("plateK14rep1_20200430_194338_941_WellG04_PointG04_0008_ChannelDAPI,WF_GFP,TRITC,WF_Cy5,DIA_Seq0683.h5", "second_pass"),
    @("titration_plate_20200403_154849_WellD01_PointD01_0004_ChannelDAPI,WF_GFP,TRITC,WF_Cy5_Seq0427.h5", "second_pass"),
    @("titration_plate_20200403_154849_WellD10_PointD10_0006_ChannelDAPI,WF_GFP,TRITC,WF_Cy5_Seq0348.h5", "second_pass"),
    @("titration_plate_20200403_154849_WellH08_PointH08_0008_ChannelDAPI,WF_GFP,TRITC,WF_Cy5_Seq0800.h5", "second_pass")
)

$startRow = 35
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("B53").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
